# Insert a new row at 173, shifting existing rows 173-201 down to 174-202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(173).Insert()

# Fill the newly inserted row 173 with the new record's data.
$ws.Range("A173").Value = 8
$ws.Range("B173").Value = "Terminal La Palmera de La Serena"
$ws.Range("C173").Value = "Coquimbo"
$ws.Range("D173").Value = 44504
$ws.Range("D173").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E173").Value = 4
$ws.Range("F173").Value = 100112032
$ws.Range("G173").Value = "Zapallo italiano"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 500
$ws.Range("K173").Value = 10000
$ws.Range("L173").Value = 10500
$ws.Range("M173").Value = 10250
$ws.Range("N173").Value = '$/caja 70 unidades'
$ws.Range("O173").Value = "Provincia de Limarí"
$ws.Range("P173").Value = 146
$ws.Range("Q173").Value = 70
$ws.Range("R173").Value = "Hortaliza"
